# Applies the "K instead of Strike#" regeneration to column G (K) in the
# crawford_kutter save_data sheet: overwrite the computed K values for
# rows 2-26 with their newly regenerated values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K (column G) values, keyed by row number, as produced by the
# regenerated std/mean calculation (s_vals).
$newK = @{
    2  = 5
    3  = 5
    4  = 4
    5  = 4
    6  = 4
    7  = 6
    8  = 2
    9  = 4
    10 = 6
    11 = 6
    12 = 8
    13 = 3
    14 = 7
    15 = 1
    16 = 2
    17 = 0
    18 = 5
    19 = 2
    20 = 1
    21 = 3
    22 = 0
    23 = 2
    24 = 5
    25 = 1
    26 = 2
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
